$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("IncomeReport")
$ws2 = $wb.Worksheets.Item("EarningsReport")

# --- Ticker + derived metrics row (shared between both sheets' A2) ---
$ws1.Range("A2").Value = "RUN"
$ws2.Range("A2").Value = "RUN"

# --- IncomeReport (sheet1) row 2 ---
$ws1.Range("B2").Value = "495.78M -> 584.58M -> 631.91M -> 609.15M -> 589.85M"
$ws1.Range("C2").Value = 19
$ws1.Range("D2").Value = "-0.42 -> -0.06 -> 0.96 -> 0.29 -> -1.12"
$ws1.Range("E2").Value = -167
$ws1.Range("F2").Value = "(683.51M) -> (716.05M) -> (636.94M) -> (823.36M) -> (949.63M)"
$ws1.Range("G2").Value = -39
$ws1.Range("H2").Value = "N/A <- N/A <- N/A <- N/A"
$ws1.Range("I2").Value = "1.90 <- 2.41 <- 2.95 <- 2.73"

# --- EarningsReport (sheet2) row 2 ---
# B2/C2 are numeric-looking text (stored as shared strings, not numbers in
# the source data), so force text entry via a leading quote, same as typing
# '-0.13 into Excel directly.
$ws2.Range("B2").Formula = "'-0.13"
$ws2.Range("C2").Formula = "'-1.12"
$ws2.Range("D2").Value = -761
$ws2.Range("E2").Value = "83, -269"
$ws2.Range("F2").Value = "7, 8"
$ws2.Range("G2").Value = -13

# --- Column width tweaks on IncomeReport (sheet1) ---
$ws1.Columns.Item(4).ColumnWidth = 19.5
$ws1.Columns.Item(6).ColumnWidth = 30.5
$ws1.Columns.Item(9).ColumnWidth = 13.8
